# Update the "Team members" list.
#
# The edit shifts each member's name/ID down one slot and inserts a new
# member ("Jay Surya") in what used to be the 2nd slot, bumping the trailing
# roll numbers by one for the members below. Net paragraph-level changes:
#
#   "Hari Babu Hasini -am.en.u4ece22122"          -> "N.Harshitha - am.en.u4ece22123"
#   "N.Harshitha - am.en.u4ece22123"               -> "Jay Surya-am.en.u4ece22124"
#   "Kothala Gnana Shanmukha -am.en.u4ece22124"    -> "Kothala Gnana Shanmukha -am.en.u4ece22125"
#   "Kudari Ritika - am.en.u4ece22125"             -> "Kudari Ritika - am.en.u4ece22126"
#
# Because the "after" text of the first paragraph equals the "before" text
# of the second paragraph, the replacements are done bottom-up (last
# paragraph first) so every search string is still unique in the document
# at the moment it is searched, and no replacement creates a premature
# duplicate of a string a later (in our processing order, but earlier in
# the document) step still needs to find.

$d = $word.ActiveDocument

$replacements = @(
    @("Kudari Ritika - am.en.u4ece22125", "Kudari Ritika - am.en.u4ece22126"),
    @("Kothala Gnana Shanmukha -am.en.u4ece22124", "Kothala Gnana Shanmukha -am.en.u4ece22125"),
    @("N.Harshitha - am.en.u4ece22123", "Jay Surya-am.en.u4ece22124"),
    @("Hari Babu Hasini -am.en.u4ece22122", "N.Harshitha - am.en.u4ece22123")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
